# Add 2022-Q1 data:
#  - insert a new "2022-Q1" sheet (fund holdings detail) right before the
#    "总计" (grand-total) summary sheet
#  - prepend a 2022-Q1 row to the "总计" summary sheet (existing rows shift
#    down by one, 2020-Q4 stays as the last row)

$wb = $excel.ActiveWorkbook

# Formatting (bold, centered, thin-bordered) in this workbook always comes
# from one of two already-registered cell styles, so the safest way to
# reproduce it on brand-new cells is to copy formatting from a cell that
# already carries it rather than re-deriving it property-by-property.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

# ---- remember the existing "总计" rows, then drop that sheet; it gets
#      rebuilt (with the new row) right after the new detail sheet so the
#      final sheet order / ids match ----
$zongji = $wb.Worksheets.Item("总计")
$summaryRows = @(
    @("2021-Q4", 1, 1.58),
    @("2021-Q3", 1, 1.39),
    @("2021-Q2", 1, 1.32),
    @("2021-Q1", 3, 3.1),
    @("2020-Q4", 3, 3.68)
)
$lastDataSheet = $wb.Worksheets.Item("2021-Q4")
[void]$zongji.Delete()

# ---- build the new "2022-Q1" fund-holdings detail sheet ----
$q1 = $wb.Worksheets.Add($null, $lastDataSheet)
$q1.Name = "2022-Q1"

$styleSrc.Range("B1").Copy($q1.Range("B1:H1"))
$styleSrc.Range("A2").Copy($q1.Range("A2"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0

$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "007497"
$q1.Range("B2").ClearFormats()

$q1.Range("C2").Value = "中庚价值灵动灵活配置混合"

$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "24.35"
$q1.Range("D2").ClearFormats()

$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "89.42"
$q1.Range("E2").ClearFormats()

$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "1.87"
$q1.Range("F2").ClearFormats()

$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.4553"
$q1.Range("G2").ClearFormats()

$q1.Range("H2").Value = 10

# ---- rebuild the "总计" summary sheet right after the new detail sheet ----
$newZongji = $wb.Worksheets.Add($null, $q1)
$newZongji.Name = "总计"

$styleSrc.Range("B1").Copy($newZongji.Range("B1:D1"))

$newZongji.Range("B1").Value = "日期"
$newZongji.Range("C1").Value = "持有数量(只)"
$newZongji.Range("D1").Value = "持有市值(亿元)"

# new 2022-Q1 row, followed by the previously existing rows shifted down by one
$allRows = @(, @("2022-Q1", 1, 0.46)) + $summaryRows

$r = 2
foreach ($row in $allRows) {
    $styleSrc.Range("A2").Copy($newZongji.Cells.Item($r, 1))
    $newZongji.Cells.Item($r, 1).Value = $r - 2
    $newZongji.Cells.Item($r, 2).Value = $row[0]
    $newZongji.Cells.Item($r, 3).Value = $row[1]
    $newZongji.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}
